$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 40.57152937025953
$ws.Range("E2").Value = -1270582.374874716
$ws.Range("D3").Value = 33.06654642467952
$ws.Range("D4").Value = 28.89170423621823
$ws.Range("E4").Value = -614977.6643250514
$ws.Range("D5").Value = 2.494511505031384
$ws.Range("E5").Value = -28540.80875100597
$ws.Range("D6").Value = 2.60139529383943
$ws.Range("E6").Value = -13463.42249435647
$ws.Range("D7").Value = 36.1356476619929
$ws.Range("E7").Value = -947681.9688098977
$ws.Range("D8").Value = 28.60110706106323
$ws.Range("E8").Value = -791374.5155412649
$ws.Range("D9").Value = 28.89170423621823
$ws.Range("D10").Value = 3.287128914228693
$ws.Range("E10").Value = -30132.46760517264
$ws.Range("D11").Value = 3.054753129461777
$ws.Range("E11").Value = -10889.96936935647
$ws.Range("D12").Value = 100.7620660457243
$ws.Range("E12").Value = -2444334.036591016
$ws.Range("D13").Value = 96.791666666667
$ws.Range("E13").Value = -2299172.33184567
$ws.Range("D14").Value = 51.05458665329401
$ws.Range("E14").Value = -646277.735017036
$ws.Range("D15").Value = 73.5113122515012
$ws.Range("E15").Value = -969922.683933706
$ws.Range("D16").Value = 1.57534246575342
$ws.Range("E16").Value = -9128.650270226879
$ws.Range("D17").Value = 147.4245115452936
$ws.Range("E17").Value = -1960898.24654109
$ws.Range("D18").Value = 123.7639553429033
$ws.Range("E18").Value = -1617897.221038187
$ws.Range("D19").Value = 67.45902991795138
$ws.Range("E19").Value = -1171732.823053524
$ws.Range("D20").Value = 71.99295011275102
$ws.Range("E20").Value = -792057.5649881973
$ws.Range("D21").Value = 6.799999999999968
$ws.Range("E21").Value = -27987.54910714286
